$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$theme = $sm.Theme
Write-Output "Theme:" $theme
try {
    $variants = $theme.ThemeVariants
    Write-Output "Variants:" $variants
    Write-Output "Variants.Count:" $variants.Count
} catch {
    Write-Output "variants failed:" $_.Exception.Message
}
